# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values in-place with freshly computed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 0
    14 = 0
    15 = 0
    16 = 2
    18 = 1
    20 = 2
    21 = 0
    22 = 3
    23 = 0
    24 = 3
    25 = 1
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 1
    39 = 2
    40 = 0
    41 = 2
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 0
    47 = 0
    48 = 1
    49 = 2
    50 = 2
    53 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
